$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "s"
$ws.Range("G2").Value = 16

$ws.Range("D2").Select()
